$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 (Line # 12) - A14 already holds 12; fill remaining columns ---
$ws.Range("B14").Value = "Detective"
$ws.Range("C14").Value = "I have a few safety concerns about your park."
$ws.Range("D14").Value = -1
$ws.Range("E14").Value = "Yes"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = "None"
$ws.Range("H14").Value = "Clicked on henchman"

# --- Row 15 (Line # 13) ---
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Henchman"
$ws.Range("C15").Value = "Your concerns are noted."
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "None"
$ws.Range("G15").Value = "None"
$ws.Range("H15").Value = "prev"

# --- Rows 16 & 17: Character column (B) filled for both before moving on ---
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("B16").Value = "Detective"
$ws.Range("B17").Value = "Bert"

# Row 16 (Line # 14) remaining columns
$ws.Range("C16").Value = "Excuse me, small child, have you seen any aliens around?"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = "No"
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = "None"
$ws.Range("H16").Value = "Clicked on Bert"

# Row 17 (Line # 15) remaining columns
$ws.Range("C17").Value = "Haven't seen one, bub. No aliens here. But..."
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = "-"
$ws.Range("F17").Value = "None"
$ws.Range("G17").Value = "None"
$ws.Range("H17").Value = "prev"

# --- Row 18 (Line # 16) ---
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Bert"
$ws.Range("C18").Value = "I lost one of my ""very special balloons"". Green one."
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "None"
$ws.Range("G18").Value = "None"
$ws.Range("H18").Value = "prev"

# --- Row 19 (Line # 17) ---
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Bert"
$ws.Range("C19").Value = "Find my balloon and I'll make it worth your while."
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "None"
$ws.Range("G19").Value = "None"
$ws.Range("H19").Value = "prev"

# --- Row 20 (Line # 18) ---
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Bert"
$ws.Range("C20").Value = "Yo, you found my ""special balloon"" yet?"
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = "Yes"
$ws.Range("F20").Value = "not has_balloon"
$ws.Range("G20").Value = "None"
$ws.Range("H20").Value = "Clicked on Bert"

# --- Row 21 (Line # 19) ---
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Detective"
$ws.Range("C21").Value = "Not yet."
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "None"
$ws.Range("G21").Value = "None"
$ws.Range("H21").Value = "prev"

# --- Column F widened (no longer shares E's bestFit width) to fit the new text ---
$ws.Columns("F").ColumnWidth = 13.25

# --- Update frozen pane scroll position and selection to match final view ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A22").Select()
